$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 164.5638641834259, 1.819545263381624, 3.958195447921753, 0.4145497810354433),
    @(3, 354.8071486473083, 18.72283649686717, 6.271195220947265, 0.6932291159734012),
    @(4, 779.4766010761261, 8.734464845723336, 6.177334833145141, 1.061880828062358),
    @(5, 296.0309174060822, 6.312961264444795, 6.502901935577393, 1.17016737817234),
    @(6, 572.5371844768524, 1.046820829642076, 5.631186389923096, 1.190985477743646),
    @(7, 1145.533498477936, 4.423988606117743, 7.072550773620605, 1.253265300936252),
    @(8, 387.1779374599457, 3.156368196832271, 5.312621593475342, 0.3734209686066884),
    @(9, 752.0191740036011, 2.210715666415394, 6.447424268722534, 0.6474194059469235),
    @(10, 1496.059204816818, 2.090992319611555, 8.130806112289429, 0.6848270952819346),
    @(11, 303.9344966888428, 2.143104431177518, 7.097133731842041, 0.6915775617205555),
    @(12, 594.8361550807953, 2.935440986060367, 7.930781650543213, 0.4163861341295256),
    @(13, 1168.329228258133, 7.47464352032362, 7.786657476425171, 0.6796401840861508),
    @(14, 432.7356279373169, 4.435914836356953, 7.687735366821289, 0.4796676226227088),
    @(15, 846.9339932918549, 10.86141788666196, 7.945489406585693, 0.5472781446604964),
    @(16, 1690.790585613251, 14.08237201106001, 8.448717880249024, 0.9075278783903197),
    @(17, 566.1129815578461, 3.41215894203568, 8.320358562469483, 0.389285599815171),
    @(18, 1098.459630298614, 6.568691610759929, 7.741063690185547, 0.6472925218204921),
    @(19, 2235.206554603576, 17.00333680274761, 8.569094944000245, 0.2319153980461384),
    @(20, 359.5666460990906, 4.054923299339243, 7.655160522460937, 0.4252345147805485),
    @(21, 704.9846969127655, 10.89732795019918, 8.03421630859375, 0.4733236186117118),
    @(22, 1387.62425866127, 18.73162652343346, 8.665713119506837, 0.2937387384025738),
    @(23, 515.4479577541351, 7.98518772219839, 8.135776281356812, 0.2557100318973322),
    @(24, 1010.814265394211, 15.95101423035066, 8.154685497283936, 0.1551566279764976),
    @(25, 2016.94050359726, 25.86675413702536, 8.412013244628906, 0.4998322299607175),
    @(26, 672.9128941535949, 6.988431497307461, 8.112998151779175, 0.1836148973457826),
    @(27, 1318.110341501236, 14.45998430083557, 7.342122220993042, 0.9033215667018971),
    @(28, 1788.452019882202, 233.5091938222065, 2.087626314163208, 1.689511123396893),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

Write-Output "Updated $($data.Count) rows (B:E) with new CI results."
